# Add two new submission rows (95 and 96) to the first worksheet
# ("八位序列号收集收集结果yd5"), matching the rows appended in the
# upstream "Add files via upload" commit.
#
# Columns: A = submitter (auto), B = submit time (auto, date/time),
#          C = serial number (required), D = QQ number (required).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("八位序列号收集收集结果yd5")

# --- Row 95 ---------------------------------------------------------
$ws.Range("A95").Value = "锡纸上的忧伤"

$ws.Range("B95").NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Range("B95").Value = 45963.5821064815

$ws.Range("C95").Value = "5efd1616"

# --- Row 96 ---------------------------------------------------------
$ws.Range("A96").Value = "Promise."

$ws.Range("B96").NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Range("B96").Value = 45963.7414467593

$ws.Range("C96").Value = "3cbc7e91"

# --- Column D (QQ numbers) -------------------------------------------
# These are all-digit strings; mark the range as Text first so Excel
# doesn't silently coerce them into numbers, then drop the leftover
# number-format override so the cells keep the sheet's plain look
# (matching the rest of column D, which carries no explicit style).
$ws.Range("D95:D96").NumberFormat = "@"
$ws.Range("D95").Value = "1039972313"
$ws.Range("D96").Value = "2926886958"
$ws.Range("D95:D96").ClearFormats()
